# Updated main GSC export data: the oldest day's row ("2025-11-19") has
# rolled off the export, so remove the first data row from the "Chart"
# sheet and let every subsequent row shift up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")
$ws.Rows(2).Delete()
